# Applies config field renames / value updates described in the commit:
# "sửa tên trường mail trong file config" (fix mail field name in config file)
#
# Changes:
#   B7  (stage)          : 2 -> 1
#   A8  (dayUpdate)      -> createdDay   (key rename, value untouched)
#   A13 (mailKhachHang)  -> mailAmber    (key rename, value untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

$ws.Range("B7").Value = 1
$ws.Range("A8").Value = "createdDay"
$ws.Range("A13").Value = "mailAmber"

# Restore selection to match the saved workbook state (cursor ends on A13)
$ws.Range("A13").Select()

$wb.Save()
